$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'69.100.82"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.24%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.745.10"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.24%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.08%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'601.76"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.10%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'167.19"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.31%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'3.743.13"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.24%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  -0.02%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.541"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +1.72%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  +3.84%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'6.38"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.51%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.18%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'37.99"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -0.25%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = "'  +2.15%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'4.370.27"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.28%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'3.750.32"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.34%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'69.111.29"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.35%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'7.35"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +1.32%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'17.36"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.50%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  -1.53%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'11.11"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +8.39%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'492.63"
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Value = "'0.727"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +0.71%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'0.0000151"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +8.61%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'84.91"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -0.39%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.37%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.88%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'10.08"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -0.51%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  -0.01%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  +1.13%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'8.13"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.82%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.69%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'31.49"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.82%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'3.891.83"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.18%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'3.679.56"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +0.35%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.04%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  +0.05%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'1.01"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +0.09%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'5.94"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +2.18%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  +4.29%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.325"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.00%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  +6.53%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'48.76"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -0.42%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'Stacks"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'1.99"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +0.08%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = "'Bittensor"
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = "'425.47"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -2.21%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = "'  +0.56%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  +0.01%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'40.13"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -1.04%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'141.92"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.21%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'2.785.55"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +1.57%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.11%  "
$ws.Range('E51').Style = 'Normal'
